$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.526385
$ws.Range("H2").Value = 1.579155
$ws.Range("I2").Value = 0.004602503788524942
$ws.Range("J2").Value = 0.004602503788524943
$ws.Range("M2").Value = 5.616015666666666
$ws.Range("N2").Value = 16.848047
$ws.Range("O2").Value = 0.2860808099623356
$ws.Range("P2").Value = 0.2860808099623357
$ws.Range("Q2").Value = 2.956186406698333
$ws.Range("R2").Value = 26.605677660285
$ws.Range("S2").Value = 0.001316688011675934
$ws.Range("T2").Value = 0.001316688011675934
$ws.Range("G3").Value = 0.526385
$ws.Range("H3").Value = 1.579155
$ws.Range("I3").Value = 0.004602503788524942
$ws.Range("J3").Value = 0.004602503788524943
$ws.Range("M3").Value = 8.435525999999999
$ws.Range("O3").Value = 0.4297071542841152
$ws.Range("P3").Value = 0.4297071542841153
$ws.Range("Q3").Value = 4.44033435351
$ws.Range("R3").Value = 39.96300918159
$ws.Range("S3").Value = 0.001977728805548912
$ws.Range("T3").Value = 0.001977728805548913
$ws.Range("G4").Value = 0.526385
$ws.Range("H4").Value = 1.579155
$ws.Range("I4").Value = 0.004602503788524942
$ws.Range("J4").Value = 0.004602503788524943
$ws.Range("M4").Value = 2.036951
$ws.Range("N4").Value = 6.110853000000001
$ws.Range("O4").Value = 0.1037626364528048
$ws.Range("P4").Value = 0.1037626364528048
$ws.Range("Q4").Value = 1.072220452135
$ws.Range("R4").Value = 9.649984069215002
$ws.Range("S4").Value = 0.0004775679273813705
$ws.Range("T4").Value = 0.0004775679273813707
$ws.Range("G5").Value = 0.526385
$ws.Range("H5").Value = 1.579155
$ws.Range("I5").Value = 0.004602503788524942
$ws.Range("J5").Value = 0.004602503788524943
$ws.Range("M5").Value = 3.542379
$ws.Range("N5").Value = 10.627137
$ws.Range("O5").Value = 0.1804493993007442
$ws.Range("P5").Value = 0.1804493993007443
$ws.Range("Q5").Value = 1.864655169915
$ws.Range("R5").Value = 16.781896529235
$ws.Range("S5").Value = 0.0008305190439187255
$ws.Range("T5").Value = 0.0008305190439187259
$ws.Range("I6").Value = 0.9930510500677584
$ws.Range("J6").Value = 0.9930510500677585
$ws.Range("M6").Value = 5.616015666666666
$ws.Range("N6").Value = 16.848047
$ws.Range("O6").Value = 0.2860808099623356
$ws.Range("P6").Value = 0.2860808099623357
$ws.Range("Q6").Value = 637.8363061182092
$ws.Range("R6").Value = 5740.526755063883
$ws.Range("S6").Value = 0.2840928487373323
$ws.Range("T6").Value = 0.2840928487373323
$ws.Range("I7").Value = 0.9930510500677584
$ws.Range("J7").Value = 0.9930510500677585
$ws.Range("M7").Value = 8.435525999999999
$ws.Range("O7").Value = 0.4297071542841152
$ws.Range("P7").Value = 0.4297071542841153
$ws.Range("Q7").Value = 958.0608501396239
$ws.Range("R7").Value = 8622.547651256615
$ws.Range("S7").Value = 0.4267211407834689
$ws.Range("T7").Value = 0.426721140783469
$ws.Range("I8").Value = 0.9930510500677584
$ws.Range("J8").Value = 0.9930510500677585
$ws.Range("M8").Value = 2.036951
$ws.Range("N8").Value = 6.110853000000001
$ws.Range("O8").Value = 0.1037626364528048
$ws.Range("P8").Value = 0.1037626364528048
$ws.Range("Q8").Value = 231.345740236324
$ws.Range("R8").Value = 2082.111662126916
$ws.Range("S8").Value = 0.1030415950872569
$ws.Range("T8").Value = 0.1030415950872569
$ws.Range("I9").Value = 0.9930510500677584
$ws.Range("J9").Value = 0.9930510500677585
$ws.Range("M9").Value = 3.542379
$ws.Range("N9").Value = 10.627137
$ws.Range("O9").Value = 0.1804493993007442
$ws.Range("P9").Value = 0.1804493993007443
$ws.Range("Q9").Value = 402.3240087525961
$ws.Range("R9").Value = 3620.916078773364
$ws.Range("S9").Value = 0.1791954654597003
$ws.Range("T9").Value = 0.1791954654597004
$ws.Range("G10").Value = 0.223143
$ws.Range("H10").Value = 0.6694290000000001
$ws.Range("I10").Value = 0.001951074789142588
$ws.Range("J10").Value = 0.001951074789142588
$ws.Range("M10").Value = 5.616015666666666
$ws.Range("N10").Value = 16.848047
$ws.Range("O10").Value = 0.2860808099623356
$ws.Range("P10").Value = 0.2860808099623357
$ws.Range("Q10").Value = 1.253174583907
$ws.Range("R10").Value = 11.278571255163
$ws.Range("S10").Value = 0.0005581650559750049
$ws.Range("T10").Value = 0.000558165055975005
$ws.Range("G11").Value = 0.223143
$ws.Range("H11").Value = 0.6694290000000001
$ws.Range("I11").Value = 0.001951074789142588
$ws.Range("J11").Value = 0.001951074789142588
$ws.Range("M11").Value = 8.435525999999999
$ws.Range("O11").Value = 0.4297071542841152
$ws.Range("P11").Value = 0.4297071542841153
$ws.Range("Q11").Value = 1.882328578218
$ws.Range("R11").Value = 16.940957203962
$ws.Range("S11").Value = 0.0008383907954379417
$ws.Range("T11").Value = 0.0008383907954379419
$ws.Range("G12").Value = 0.223143
$ws.Range("H12").Value = 0.6694290000000001
$ws.Range("I12").Value = 0.001951074789142588
$ws.Range("J12").Value = 0.001951074789142588
$ws.Range("M12").Value = 2.036951
$ws.Range("N12").Value = 6.110853000000001
$ws.Range("O12").Value = 0.1037626364528048
$ws.Range("P12").Value = 0.1037626364528048
$ws.Range("Q12").Value = 0.4545313569930001
$ws.Range("R12").Value = 4.090782212937
$ws.Range("S12").Value = 0.0002024486640380352
$ws.Range("T12").Value = 0.0002024486640380353
$ws.Range("G13").Value = 0.223143
$ws.Range("H13").Value = 0.6694290000000001
$ws.Range("I13").Value = 0.001951074789142588
$ws.Range("J13").Value = 0.001951074789142588
$ws.Range("M13").Value = 3.542379
$ws.Range("N13").Value = 10.627137
$ws.Range("O13").Value = 0.1804493993007442
$ws.Range("P13").Value = 0.1804493993007443
$ws.Range("Q13").Value = 0.7904570771970001
$ws.Range("R13").Value = 7.114113694773001
$ws.Range("S13").Value = 0.0003520702736916063
$ws.Range("T13").Value = 0.0003520702736916064
$ws.Range("G14").Value = 0.04521833333333333
$ws.Range("H14").Value = 0.135655
$ws.Range("I14").Value = 0.0003953713545740292
$ws.Range("J14").Value = 0.0003953713545740293
$ws.Range("M14").Value = 5.616015666666666
$ws.Range("N14").Value = 16.848047
$ws.Range("O14").Value = 0.2860808099623356
$ws.Range("P14").Value = 0.2860808099623357
$ws.Range("Q14").Value = 0.2539468684205555
$ws.Range("R14").Value = 2.285521815785
$ws.Range("S14").Value = 0.0001131081573524441
$ws.Range("T14").Value = 0.0001131081573524441
$ws.Range("G15").Value = 0.04521833333333333
$ws.Range("H15").Value = 0.135655
$ws.Range("I15").Value = 0.0003953713545740292
$ws.Range("J15").Value = 0.0003953713545740293
$ws.Range("M15").Value = 8.435525999999999
$ws.Range("O15").Value = 0.4297071542841152
$ws.Range("P15").Value = 0.4297071542841153
$ws.Range("Q15").Value = 0.38144042651
$ws.Range("R15").Value = 3.43296383859
$ws.Range("S15").Value = 0.000169893899659462
$ws.Range("T15").Value = 0.000169893899659462
$ws.Range("G16").Value = 0.04521833333333333
$ws.Range("H16").Value = 0.135655
$ws.Range("I16").Value = 0.0003953713545740292
$ws.Range("J16").Value = 0.0003953713545740293
$ws.Range("M16").Value = 2.036951
$ws.Range("N16").Value = 6.110853000000001
$ws.Range("O16").Value = 0.1037626364528048
$ws.Range("P16").Value = 0.1037626364528048
$ws.Range("Q16").Value = 0.09210752930166667
$ws.Range("R16").Value = 0.828967763715
$ws.Range("S16").Value = 0.00004102477412851799
$ws.Range("T16").Value = 0.00004102477412851799
$ws.Range("G17").Value = 0.04521833333333333
$ws.Range("H17").Value = 0.135655
$ws.Range("I17").Value = 0.0003953713545740292
$ws.Range("J17").Value = 0.0003953713545740293
$ws.Range("M17").Value = 3.542379
$ws.Range("N17").Value = 10.627137
$ws.Range("O17").Value = 0.1804493993007442
$ws.Range("P17").Value = 0.1804493993007443
$ws.Range("Q17").Value = 0.160180474415
$ws.Range("R17").Value = 1.441624269735
$ws.Range("S17").Value = 0.00007134452343360513
$ws.Range("T17").Value = 0.00007134452343360516
